# edit.ps1 - Apply quarterly financial update to DOV_QTR_FIN workbook
# Inserts two new columns (new quarters: 2018-11-30 and 2018-09-30 period-ends,
# serials 43465 and 43373) before column D, shifting the existing quarterly
# history from D:K to F:M, fills in the two new columns with the latest
# reported figures, and applies a handful of restatements that the data
# provider issued for the two quarters that land in columns H and I (and J
# for one row) after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns before column D; this shifts old D:K -> F:M
#    and preserves the values/number formats that were already there.
$ws.Columns("D:E").Insert()

# 2) The newly inserted D:E columns inherit column C's style; restyle them to
#    match the rest of the data columns (style index 3 = number style, except
#    the three "Period Ending" rows which use the date style index 2).
$ws.Range("D5:E102").Style = $ws.Range("F5").Style

$dateRows = @(7,38,80)
foreach ($r in $dateRows) {
    $ws.Range("D" + $r + ":E" + $r).Style = $ws.Range("F" + $r).Style
}

# 3) Populate the two new columns (D, E) with the newest reported quarters.
$newData = @(
    @(7, 43465, 43373),
    @(8, 1809000, 1747400),
    @(9, 1155000, 1097300),
    @(10, 654000, 650100),
    @(11, $null, $null),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 22100, 24200),
    @(15, 0, 0),
    @(16, $null, $null),
    @(17, 1590200, 1527300),
    @(18, 218800, 220100),
    @(19, $null, $null),
    @(20, -100, 4100),
    @(21, 295300, 292300),
    @(22, 32000, 31200),
    @(23, 186700, 193000),
    @(24, 31600, 35700),
    @(25, 0, 0),
    @(26, 155100, 157300),
    @(27, 155100, 157300),
    @(28, 0, 0),
    @(29, -13500, 0),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 100, -4100),
    @(33, 141600, 157300),
    @(34, 0, 0),
    @(35, 141600, 157300),
    @(38, 43465, 43373),
    @(39, $null, $null),
    @(40, $null, $null),
    @(41, 396200, 209300),
    @(42, 0, 0),
    @(43, 1241200, 1300000),
    @(44, 748800, 804000),
    @(45, 117500, 145800),
    @(46, 2503800, 2459000),
    @(47, 0, 0),
    @(48, 806500, 806700),
    @(49, 4811600, 4915400),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 243900, 280500),
    @(53, 0, 0),
    @(54, 8365800, 8461800),
    @(55, $null, $null),
    @(56, $null, $null),
    @(57, 969500, 948600),
    @(58, 220300, 298700),
    @(59, 637600, 643000),
    @(60, 1827400, 1890300),
    @(61, 2943700, 2981900),
    @(62, 826000, 841600),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 5597100, 5713800),
    @(67, $null, $null),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 7815500, 7744400),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 2768700, 2747900),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 141600, 157300),
    @(82, $null, $null),
    @(83, 76600, 68100),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 364200, 240400),
    @(90, $null, $null),
    @(91, -36400, -38200),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -35400, -35900),
    @(95, $null, $null),
    @(96, -70400, -70800),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, -148500, -232500),
    @(101, 6700, -5500),
    @(102, 186900, -33500)
)

foreach ($item in $newData) {
    $r = $item[0]
    $dVal = $item[1]
    $eVal = $item[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($r, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($r, 5).Value = $eVal
    }
}

# 4) A handful of historical quarters (now in columns H, I, and one in J)
#    were restated by the data provider along with this update; apply those
#    corrected figures on top of the shifted values.
$restated = @(
    @("H8", 1752500),
    @("I8", 1747800),
    @("H9", 1092100),
    @("I9", 1096700),
    @("H10", 660400),
    @("I10", 651100),
    @("H14", -73100),
    @("H17", 1454700),
    @("I17", 1508600),
    @("H18", 297800),
    @("I18", 239200),
    @("H20", 2700),
    @("I20", 3000),
    @("H21", 371600),
    @("I21", 314700),
    @("I22", 35400),
    @("H23", 264100),
    @("I23", 206800),
    @("H24", 25400),
    @("I24", 47300),
    @("H26", 238700),
    @("I26", 159500),
    @("H27", 238700),
    @("I27", 159500),
    @("H29", 57700),
    @("I29", 19500),
    @("H32", -2700),
    @("I32", -3000),
    @("H89", 323800),
    @("I89", 273500),
    @("H91", -39700),
    @("I91", -51400),
    @("J91", -36700),
    @("H100", -86900)
)

foreach ($item in $restated) {
    $addr = $item[0]
    $val = $item[1]
    $ws.Range($addr).Value = $val
}
